$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up floating-point noise in B2, B4, B7 (values are numerically the same,
# just re-written with a cleaner decimal representation)
$ws.Range("B2").Value = 55960.978
$ws.Range("B4").Value = 89299.253
$ws.Range("B7").Value = 222082.585

# Update B6 with the new value (+72.7 change)
$ws.Range("B6").Value = 34342.29500000001
